$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 00:22"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 394278
$ws.Range("C4").Value = 27274
$ws.Range("D4").Value = 21650
$ws.Range("E4").Value = 359911
$ws.Range("G4").Value = 1846
$ws.Range("H4").Value = 12717

# Row 5 - España
$ws.Range("B5").Value = 141942
$ws.Range("C5").Value = 5267
$ws.Range("E5").Value = 84689
$ws.Range("G5").Value = 704
$ws.Range("H5").Value = 14045

# Row 8 - Alemania
$ws.Range("B8").Value = 107663
$ws.Range("C8").Value = 4288
$ws.Range("E8").Value = 69566
$ws.Range("G8").Value = 206
$ws.Range("H8").Value = 2016

# Row 9 - China
$ws.Range("C9").Value = 0

# Row 17 - Brasil
$ws.Range("B17").Value = 14018
$ws.Range("C17").Value = 1835
$ws.Range("E17").Value = 13205
$ws.Range("G17").Value = 122
$ws.Range("H17").Value = 686

# Row 130 - Madagascar
$ws.Range("D130").Value = 7
$ws.Range("E130").Value = 81

# Row 139 - Jamaica
$ws.Range("B139").Value = 63
$ws.Range("C139").Value = 4
$ws.Range("E139").Value = 52

# Row 153 - Guyana
$ws.Range("F153").Value = 4
